# Update curriculum vitae workbook:
#  - Academic History sheet: refresh master's-thesis link/title to the
#    Portuguese title + new https link, fix the bachelor's monograph
#    hyperlink, turn the PhD thesis link into a real hyperlink, and move
#    the IDP/MBA "academic" row out to the Training sheet.
#  - Training sheet: receives the IDP / MBA in Big Data row that used to
#    live on Academic History.
#  - Make "Academic History" the active/selected tab.

$wb = $excel.ActiveWorkbook

$wsWork     = $wb.Worksheets.Item("Work Experience")
$wsAcademic = $wb.Worksheets.Item("Academic History")
$wsTraining = $wb.Worksheets.Item("Training")

# ---------------------------------------------------------------------
# 1. Academic History row 3 (M.S. Degree @ ENCE/IBGE): swap the English
#    dissertation title for the Portuguese one, and point the link at the
#    new https ence.ibge.gov.br URL (no "www").
# ---------------------------------------------------------------------
$newDissertationTitle = "Juventude e direito à cidade: mobilizações estudantis nas instituições de ensino superior na Região Metropolitana do Rio de Janeiro em 2016"
$newDissertationUrl   = "https://ence.ibge.gov.br/images/ence/doc/mestrado/dissertacoes/2018/Dissertacao_CauanBraga_2018.pdf"

$wsAcademic.Range("G3").Value = $newDissertationTitle
$wsAcademic.Range("H3").Value = $newDissertationUrl

foreach ($hl in $wsAcademic.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$H$3') {
        $hl.Address = $newDissertationUrl
    }
}

# ---------------------------------------------------------------------
# 2. Academic History row 4 (PhD @ ENCE/IBGE): the thesis link cell was
#    plain text before; make it a proper hyperlink like the rows above it.
# ---------------------------------------------------------------------
$thesisUrl = "https://ence.ibge.gov.br/images/ence/pos_graduacao/seminarios_de_dissertacao/2022/Tese_Cauan_Braga_da_Silva_Cardoso_1.pdf"
[void]$wsAcademic.Hyperlinks.Add($wsAcademic.Range("H4"), $thesisUrl)
$wsAcademic.Range("H4").Style = "Hiperlink"

# ---------------------------------------------------------------------
# 3. Move the IDP / MBA row (row 5) from Academic History to a new last
#    row on the Training sheet, then remove it from Academic History.
# ---------------------------------------------------------------------
$idpInstitution = $wsAcademic.Range("A5").Value2
$idpAbbrev      = $wsAcademic.Range("B5").Value2
$idpTitle       = $wsAcademic.Range("C5").Value2
$idpStart       = $wsAcademic.Range("D5").Value2
$idpEnd         = $wsAcademic.Range("E5").Value2

$trainingRow = $wsTraining.Cells.Item($wsTraining.UsedRange.Rows.Count + 1, 1).Row

$wsTraining.Range("A" + $trainingRow).Value = $idpInstitution
$wsTraining.Range("B" + $trainingRow).Value = $idpAbbrev
$wsTraining.Range("C" + $trainingRow).Value = $idpTitle
$wsTraining.Range("D" + $trainingRow).Value = $idpStart
$wsTraining.Range("E" + $trainingRow).Value = $idpEnd
$wsTraining.Range("F" + $trainingRow).Value = 457

# Formats: A/B like the "Institution"/"Abbreviation" columns used on the
# Work Experience sheet; D/E like the existing Training date columns.
$wsWork.Range("A3").Copy()
$wsTraining.Range("A" + $trainingRow).PasteSpecial(-4122)
$wsWork.Range("B3").Copy()
$wsTraining.Range("B" + $trainingRow).PasteSpecial(-4122)
$wsTraining.Range("D4").Copy()
$wsTraining.Range("D" + $trainingRow).PasteSpecial(-4122)
$wsTraining.Range("E4").Copy()
$wsTraining.Range("E" + $trainingRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Now drop the row from Academic History.
[void]$wsAcademic.Rows.Item(5).Delete()

# ---------------------------------------------------------------------
# 4. Tab / selection bookkeeping: Academic History becomes the active
#    sheet, with its own selection parked on the monograph title cell;
#    Training's selection moves onto the freshly-added row.
# ---------------------------------------------------------------------
[void]$wsTraining.Activate()
[void]$wsTraining.Range("A" + $trainingRow).Select()

[void]$wsAcademic.Activate()
[void]$wsAcademic.Range("G2").Select()
